# Refresh the cryptocurrency price/volume table on Sheet1 (columns D=Price,
# E=Volume(1h)) with freshly scraped values, mirroring the GitHub Actions
# "Updated cryptos list" bot commit.
#
# Several Price values are plain digit-and-dot strings (e.g. "0.9995",
# "44.90") that Excel would otherwise auto-convert to numbers (losing
# formatting such as trailing zeros) when assigned to a General-formatted
# cell, so those specific cells are switched to Text format ("@") first to
# keep them as literal text, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.037.96"
$ws.Range("E2").Value = "  -1.37%  "

$ws.Range("D3").Value = "1.840.02"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.59%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "278.61"
$ws.Range("E5").Value = "  -2.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  -0.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5097"
$ws.Range("E7").Value = "  -1.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3497"
$ws.Range("E8").Value = "  -4.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.90"
$ws.Range("E9").Value = "  +1.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06816"
$ws.Range("E10").Value = "  -3.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.90"
$ws.Range("E11").Value = "  -3.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8092"
$ws.Range("E12").Value = "  -6.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07770"
$ws.Range("E13").Value = "  -2.78%  "

$ws.Range("D14").Value = "1.835.37"
$ws.Range("E14").Value = "  -0.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.080"
$ws.Range("E15").Value = "  -1.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.32"
$ws.Range("E16").Value = "  -1.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9994"

$ws.Range("E18").Value = "  -1.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008046"
$ws.Range("E19").Value = "  -3.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9994"
$ws.Range("E20").Value = "  -0.33%  "

$ws.Range("D21").Value = "26.093.71"
$ws.Range("E21").Value = "  -1.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.775"
$ws.Range("E22").Value = "  -2.19%  "

$ws.Range("E23").Value = "  -3.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.211"
$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.363"
$ws.Range("E25").Value = "  +7.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.56"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.661"
$ws.Range("E27").Value = "  -3.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.19"
$ws.Range("E28").Value = "  -2.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "109.65"
$ws.Range("E29").Value = "  -1.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.357"
$ws.Range("E30").Value = "  -4.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.288"
$ws.Range("E31").Value = "  -4.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08804"
$ws.Range("E32").Value = "  -1.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04860"
$ws.Range("E33").Value = "  -0.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.163"
$ws.Range("E34").Value = "  +2.24%  "

$ws.Range("E35").Value = "  -4.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.865"
$ws.Range("E36").Value = "  +0.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.214"
$ws.Range("E37").Value = "  +2.85%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.373"
$ws.Range("E38").Value = "  -6.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01851"
$ws.Range("E39").Value = "  -2.55%  "

$ws.Range("E40").Value = "  -9.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9527"
$ws.Range("E41").Value = "  -7.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "117.01"
$ws.Range("E42").Value = "  +3.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.259"
$ws.Range("E43").Value = "  -1.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.000"
$ws.Range("E44").Value = "  -5.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9983"
$ws.Range("E45").Value = "  -0.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4521"
$ws.Range("E46").Value = "  -9.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1361"
$ws.Range("E47").Value = "  -6.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.284"
$ws.Range("E48").Value = "  -4.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.13"
$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05914"
$ws.Range("E50").Value = "  -0.91%  "

$ws.Range("E51").Value = "  -4.18%  "
